$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("!!Compartment")
$ws.Unprotect()
$ws.Range("A1").Value = "!!!ObjTables schema='SBtab' objTablesVersion='0.0.8' date='2020-03-09 23:59:06'"
$ws.Range("A2").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Compartment' name='Compartment' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Compound")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Compound' name='Compound' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Definition")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Definition' name='Definition' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Enzyme")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Enzyme' name='Enzyme' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!FbcObjective")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='FbcObjective' name='FbcObjective' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Gene")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Gene' name='Gene' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Layout")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Layout' name='Layout' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Measurement")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Measurement' name='Measurement' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!PbConfig")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='PbConfig' name='PbConfig' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Position")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Position' name='Position' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Protein")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Protein' name='Protein' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Quantity")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Quantity' name='Quantity' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!QuantityInfo")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='QuantityInfo' name='QuantityInfo' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!QuantityMatrix")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='QuantityMatrix' name='QuantityMatrix' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Reaction")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Reaction' name='Reaction' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!ReactionStoichiometry")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='ReactionStoichiometry' name='ReactionStoichiometry' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Regulator")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Regulator' name='Regulator' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Relation")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Relation' name='Relation' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!Relationship")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='Relationship' name='Relationship' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!SparseMatrix")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='SparseMatrix' name='SparseMatrix' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!SparseMatrixColumn")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='SparseMatrixColumn' name='SparseMatrixColumn' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!SparseMatrixOrdered")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='SparseMatrixOrdered' name='SparseMatrixOrdered' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!SparseMatrixRow")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='SparseMatrixRow' name='SparseMatrixRow' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!StoichiometricMatrix")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='StoichiometricMatrix' name='StoichiometricMatrix' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!rxnconContingencyList")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='rxnconContingencyList' name='rxnconContingencyList' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()

$ws = $wb.Worksheets.Item("!!rxnconReactionList")
$ws.Unprotect()
$ws.Range("A1").Value = "!!ObjTables schema='SBtab' type='Data' tableFormat='row' id='rxnconReactionList' name='rxnconReactionList' date='2020-03-09 23:59:06' objTablesVersion='0.0.8'"
$ws.Protect()
